# This edit re-orders the 11 data rows (rows 2-12) of the sheet. Rows 6 and 10
# keep their original contents; the rest are rewritten in place with the data
# from another row of the original table (a straight row permutation), per the
# mapping: new row -> source (original) row.
#   2<-8   3<-4   4<-11   5<-9   6<-6(unchanged)   7<-5   8<-3   9<-2   10<-10(unchanged)   11<-12   12<-7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the original values (columns A..R) for rows 2..12 before mutating anything,
# so that the permutation reads consistently from the pre-edit state.
# NOTE: use Value() (explicit getter call) rather than the bare .Value property,
# which in this shim returns property-metadata instead of invoking the getter.
$original = @{}
for ($r = 2; $r -le 12; $r++) {
    $row = @{}
    for ($c = 1; $c -le 18; $c++) {
        $row[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $original[$r] = $row
}

# Mapping: destination row -> source row (in terms of the original, pre-edit data)
$mapping = @{
    2  = 8
    3  = 4
    4  = 11
    5  = 9
    6  = 6
    7  = 5
    8  = 3
    9  = 2
    10 = 10
    11 = 12
    12 = 7
}

foreach ($destRow in ($mapping.Keys | Sort-Object)) {
    $srcRow = $mapping[$destRow]
    if ($srcRow -eq $destRow) {
        continue
    }
    $srcData = $original[$srcRow]
    for ($c = 1; $c -le 18; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $srcData[$c]
    }
}
